$d = $word.ActiveDocument

# --- Change 1: remove the _GoBack bookmark between "," and " Eric Kerrigan" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 2: flesh out the last (empty) bulleted paragraph and append the
#     rest of the Gantt-chart / Newtonian-model discussion bullets ---

# The document currently ends with a single empty ListParagraph bullet
# (ilvl 0, numId 2). Fill it in, then append the remaining bullets below it,
# each inheriting style/numbering from the paragraph before it.
$last = $d.Paragraphs.Last
$last.Range.InsertAfter("Gantt chart discussion")

$items = @(
    @{ ilvl = 2; text = "Workflow"; bold = $false },
    @{ ilvl = 3; text = "Setup with 3 agents, centralized computation, move from A to B, constant height, no obstacles"; bold = $false },
    @{ ilvl = 3; text = "Individual computation, move from A to B, constant height, no obstacles"; bold = $false },
    @{ ilvl = 3; text = "Individual computation, move from A to B, with obstacles"; bold = $false },
    @{ ilvl = 3; text = "Implement takeoff, cruising and landing modes"; bold = $false },
    @{ ilvl = 1; text = "Newtonian model equations are being transcribed for use inside ICLOCS"; bold = $false },
    @{ ilvl = 2; text = "Unsure about how to implement differential algebraic equations"; bold = $false },
    @{ ilvl = 2; text = [string]::Concat("Read Betts", [char]0x2019, " chapter on DAEs"); bold = $true },
    @{ ilvl = 2; text = "Consult Yuanbo about DAEs in ICLOCS"; bold = $true }
)

foreach ($item in $items) {
    $prev = $d.Paragraphs.Last
    $prev.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Last
    $start = $cur.Range.Start
    $cur.Range.InsertAfter($item.text)
    $end = $cur.Range.End
    $cur.Range.ListFormat.ListLevelNumber = $item.ilvl
    if ($item.bold) {
        $rr = $d.Range($start, $end)
        $rr.Font.Bold = 1
    }
}

# Re-add the _GoBack bookmark immediately before the text of the final bullet
$final = $d.Paragraphs.Last
$bmPos = $final.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

Write-Output "done"
